# Scheduled runner update: refresh market-price derived profit columns (H,I,J,K,L,M,N)
# across several crafting-class profit sheets. Values mirror the latest source data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 44986.625
$ws.Range("I21").Value = 54982.168
$ws.Range("K21").Value = 54982.168
$ws.Range("M21").Value = -54514.168
$ws.Range("H23").Value = 44986.625
$ws.Range("I23").Value = 54982.168
$ws.Range("K23").Value = 54982.168
$ws.Range("M23").Value = -54748.168
$ws.Range("H29").Value = 657.875
$ws.Range("I29").Value = 657.875
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1973.625
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1692.625
$ws.Range("N29").ClearContents()
$ws.Range("H33").Value = 164.16667
$ws.Range("I33").Value = 97
$ws.Range("K33").Value = 97
$ws.Range("M33").Value = 132
$ws.Range("H113").Value = 2173.125
$ws.Range("I113").Value = 1712.5
$ws.Range("J113").Value = 2633.75
$ws.Range("K113").Value = 1712.5
$ws.Range("L113").Value = 2633.75
$ws.Range("M113").Value = 1541.5
$ws.Range("N113").Value = -9141.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2280.5
$ws.Range("I45").Value = 2052.75
$ws.Range("J45").Value = 2432.3333
$ws.Range("K45").Value = 2052.75
$ws.Range("L45").Value = 2432.3333
$ws.Range("M45").Value = -1675.75
$ws.Range("N45").Value = -3186.3333
$ws.Range("H122").Value = 1513
$ws.Range("I122").Value = 992.75
$ws.Range("J122").Value = 1833.1538
$ws.Range("K122").Value = 2978.25
$ws.Range("L122").Value = 5499.4614
$ws.Range("M122").Value = -528.25
$ws.Range("N122").Value = -10399.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 125.933334
$ws.Range("I64").Value = 145.33333
$ws.Range("J64").Value = 113
$ws.Range("K64").Value = 145.33333
$ws.Range("L64").Value = 113
$ws.Range("M64").Value = 79.66667000000001
$ws.Range("N64").Value = -563
$ws.Range("H67").Value = 125.933334
$ws.Range("I67").Value = 145.33333
$ws.Range("J67").Value = 113
$ws.Range("K67").Value = 145.33333
$ws.Range("L67").Value = 113
$ws.Range("M67").Value = 634.6666700000001
$ws.Range("N67").Value = -1673

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 977.1875
$ws.Range("I122").Value = 912.8
$ws.Range("J122").Value = 1084.5
$ws.Range("K122").Value = 2738.4
$ws.Range("L122").Value = 3253.5
$ws.Range("M122").Value = -288.3999999999996
$ws.Range("N122").Value = -8153.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 754.6087
$ws.Range("J131").Value = 977.875
$ws.Range("L131").Value = 2933.625
$ws.Range("N131").Value = -13013.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 29.545454
$ws.Range("I2").Value = 16.866667
$ws.Range("K2").Value = 16.866667
$ws.Range("M2").Value = 96.13333299999999
$ws.Range("H122").Value = 1154.3914
$ws.Range("I122").Value = 1079.3572
$ws.Range("J122").Value = 1271.1111
$ws.Range("K122").Value = 3238.0716
$ws.Range("L122").Value = 3813.3333
$ws.Range("M122").Value = -788.0715999999998
$ws.Range("N122").Value = -8713.3333
$ws.Range("H126").Value = 1662.6154
$ws.Range("I126").Value = 1440.4
$ws.Range("J126").Value = 1801.5
$ws.Range("K126").Value = 4321.200000000001
$ws.Range("L126").Value = 5404.5
$ws.Range("M126").Value = -1851.200000000001
$ws.Range("N126").Value = -10344.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1840.9131
$ws.Range("I7").Value = 1808.8
$ws.Range("J7").Value = 1901.125
$ws.Range("K7").Value = 1808.8
$ws.Range("L7").Value = 1901.125
$ws.Range("M7").Value = -1696.8
$ws.Range("N7").Value = -2125.125
$ws.Range("H46").Value = 1974.4
$ws.Range("I46").Value = 704.25
$ws.Range("J46").Value = 2436.2727
$ws.Range("K46").Value = 704.25
$ws.Range("L46").Value = 2436.2727
$ws.Range("M46").Value = -516.25
$ws.Range("N46").Value = -2812.2727
$ws.Range("H61").Value = 1525.8
$ws.Range("I61").Value = 1150.2222
$ws.Range("J61").Value = 2089.1667
$ws.Range("K61").Value = 1150.2222
$ws.Range("L61").Value = 2089.1667
$ws.Range("M61").Value = -948.2221999999999
$ws.Range("N61").Value = -2493.1667
$ws.Range("H82").Value = 1359
$ws.Range("I82").Value = 1322
$ws.Range("J82").Value = 1381.2
$ws.Range("K82").Value = 1322
$ws.Range("L82").Value = 1381.2
$ws.Range("M82").Value = -961
$ws.Range("N82").Value = -2103.2
$ws.Range("H85").Value = 1359
$ws.Range("I85").Value = 1322
$ws.Range("J85").Value = 1381.2
$ws.Range("K85").Value = 1322
$ws.Range("L85").Value = 1381.2
$ws.Range("M85").Value = -74
$ws.Range("N85").Value = -3877.2
$ws.Range("H113").Value = 1525.8
$ws.Range("I113").Value = 1150.2222
$ws.Range("J113").Value = 2089.1667
$ws.Range("K113").Value = 1150.2222
$ws.Range("L113").Value = 2089.1667
$ws.Range("M113").Value = 1019.7778
$ws.Range("N113").Value = -6429.1667
$ws.Range("H126").Value = 1840.9131
$ws.Range("I126").Value = 1808.8
$ws.Range("J126").Value = 1901.125
$ws.Range("K126").Value = 5426.4
$ws.Range("L126").Value = 5703.375
$ws.Range("M126").Value = -2956.4
$ws.Range("N126").Value = -10643.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2900.2856
$ws.Range("I62").Value = 2958.8
$ws.Range("J62").Value = 2867.7778
$ws.Range("K62").Value = 2958.8
$ws.Range("L62").Value = 2867.7778
$ws.Range("M62").Value = -2334.8
$ws.Range("N62").Value = -4115.7778
$ws.Range("H65").Value = 2900.2856
$ws.Range("I65").Value = 2958.8
$ws.Range("J65").Value = 2867.7778
$ws.Range("K65").Value = 14794
$ws.Range("L65").Value = 14338.889
$ws.Range("M65").Value = -11674
$ws.Range("N65").Value = -20578.889
$ws.Range("H100").Value = 891.82355
$ws.Range("I100").Value = 797
$ws.Range("J100").Value = 1334.3334
$ws.Range("K100").Value = 1594
$ws.Range("L100").Value = 2668.6668
$ws.Range("M100").Value = -1053
$ws.Range("N100").Value = -3750.6668
$ws.Range("H110").Value = 17000
$ws.Range("J110").Value = 17000
$ws.Range("L110").Value = 17000
$ws.Range("N110").Value = -25180
$ws.Range("H122").Value = 9092191
$ws.Range("I122").Value = 16667492
$ws.Range("J122").Value = 1831
$ws.Range("K122").Value = 50002476
$ws.Range("L122").Value = 5493
$ws.Range("M122").Value = -50000026
$ws.Range("N122").Value = -10393
$ws.Range("H126").Value = 923.3889
$ws.Range("I126").Value = 867.2
$ws.Range("J126").Value = 945
$ws.Range("K126").Value = 2601.6
$ws.Range("L126").Value = 2835
$ws.Range("M126").Value = -131.6000000000004
$ws.Range("N126").Value = -7775
